$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on price column cells whose new values would otherwise
# be auto-converted to numeric/date types by Excel, so they remain plain text like the rest of the sheet.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.933.06"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.630.02"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "211.75"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -1.71%  "
$ws.Range("D9").Value = "0.257"
$ws.Range("E9").Value = "  -2.10%  "
$ws.Range("D10").Value = "0.0614"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "1.860.97"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "1.624.94"
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").Value = "0.564"
$ws.Range("E15").Value = "  -2.46%  "
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").Value = "27.922.22"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "230.78"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").Value = "10.20"
$ws.Range("E23").Value = "  -5.60%  "
$ws.Range("E24").Value = "  -1.94%  "
$ws.Range("D25").Value = "155.00"
$ws.Range("E25").Value = "  +2.26%  "
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "15.56"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("D31").Value = "0.0483"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("D34").Value = "1.403.14"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").Value = "0.993"
$ws.Range("E36").Value = "  +7.49%  "
$ws.Range("D37").Value = "2.35"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("E38").Value = "  +1.60%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "0.864"
$ws.Range("E40").Value = "  -3.16%  "
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("D45").Value = "66.03"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("D47").Value = "1.770.65"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").Value = "88.20"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.101"
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.0504"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.54"
$ws.Range("E51").Value = "  -1.26%  "

# Reset style index back to default (Normal) so no stray style attribute is left on edited cells,
# while keeping the values as text (the "@" number format already forced text interpretation).
$ws.Range("D2:D51").Style = "Normal"
